$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new tag values into the "cq:tags@[]" column (F) for the rows whose
# other columns (source/target) are already filled in, matching the new
# shared strings added to the sheet.
$ws.Range("F2").Value = "properties:orientation/square"
$ws.Range("F4").Value = "properties:orientation/landscape,properties:style/color"
$ws.Range("F5").Value = "properties:orientation/landscape,properties:style/monochrome"
$ws.Range("F6").Value = "properties:orientation/landscape,properties:style/color"

# Wrap the text in those tag cells (adds a new wrapText cell style).
$ws.Range("F2").WrapText = $true
$ws.Range("F4").WrapText = $true
$ws.Range("F5").WrapText = $true
$ws.Range("F6").WrapText = $true

# The wrapped, multi-tag cells need a taller row to show their contents.
$ws.Rows.Item(4).RowHeight = 28.5
$ws.Rows.Item(5).RowHeight = 28.5
$ws.Rows.Item(6).RowHeight = 28.5

# Column F can narrow slightly now that its text wraps instead of
# overflowing. (26.1666... chars round-trips to a stored width of 27.)
$ws.Columns.Item(6).ColumnWidth = 26.16666666666667

# Scroll the view over one column and move the active selection/cursor.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("F8").Select() | Out-Null
